# Update Sheets via scheduled runner - apply updated market price data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1598.407
$ws.Range("I15").Value = 1598.407
$ws.Range("K15").Value = 4795.221
$ws.Range("M15").Value = -4626.221

# row 129
$ws.Range("H129").Value = 1582.8733
$ws.Range("I129").Value = 734
$ws.Range("J129").Value = 2044.2174
$ws.Range("K129").Value = 2202
$ws.Range("L129").Value = 6132.6522
$ws.Range("M129").Value = 2798
$ws.Range("N129").Value = -16132.6522

# row 135
$ws.Range("H135").Value = 35714940
$ws.Range("I135").Value = 341.0909
$ws.Range("K135").Value = 3069.8181
$ws.Range("M135").Value = -534.8181

$ws = $wb.Worksheets.Item("ARM")
# row 44
$ws.Range("H44").Value = 13171
$ws.Range("J44").Value = 13171
$ws.Range("L44").Value = 13171
$ws.Range("N44").Value = -14147

# row 55
$ws.Range("H55").Value = 7000
$ws.Range("J55").Value = 6857.143
$ws.Range("L55").Value = 6857.143
$ws.Range("N55").Value = -7487.143

# row 61
$ws.Range("H61").Value = 2424.5
$ws.Range("I61").Value = 2335.3333
$ws.Range("J61").Value = 2692
$ws.Range("K61").Value = 2335.3333
$ws.Range("L61").Value = 2692
$ws.Range("M61").Value = -2123.3333
$ws.Range("N61").Value = -3116

# row 74
$ws.Range("H74").Value = 32392.06
$ws.Range("I74").Value = 54755.26
$ws.Range("J74").Value = 2042
$ws.Range("K74").Value = 54755.26
$ws.Range("L74").Value = 2042
$ws.Range("M74").Value = -53881.26
$ws.Range("N74").Value = -3790

# row 77
$ws.Range("H77").Value = 32392.06
$ws.Range("I77").Value = 54755.26
$ws.Range("J77").Value = 2042
$ws.Range("K77").Value = 273776.3
$ws.Range("L77").Value = 10210
$ws.Range("M77").Value = -269408.3
$ws.Range("N77").Value = -18946

# row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# row 132
$ws.Range("H132").Value = 18932342
$ws.Range("I132").Value = 24131120
$ws.Range("J132").Value = 3336002.8
$ws.Range("K132").Value = 72393360
$ws.Range("L132").Value = 10008008.4
$ws.Range("M132").Value = -72390830
$ws.Range("N132").Value = -10013068.4

# row 136
$ws.Range("H136").Value = 2424.5
$ws.Range("I136").Value = 2335.3333
$ws.Range("J136").Value = 2692
$ws.Range("K136").Value = 7005.999899999999
$ws.Range("L136").Value = 8076
$ws.Range("M136").Value = -4455.999899999999
$ws.Range("N136").Value = -13176

$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 13545440
$ws.Range("I134").Value = 26320644
$ws.Range("J134").Value = 60503
$ws.Range("K134").Value = 78961932
$ws.Range("L134").Value = 181509
$ws.Range("M134").Value = -78959397
$ws.Range("N134").Value = -186579

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1815.6666
$ws.Range("I31").Value = 883.12
$ws.Range("J31").Value = 2829.3044
$ws.Range("K31").Value = 883.12
$ws.Range("L31").Value = 2829.3044
$ws.Range("M31").Value = -588.12
$ws.Range("N31").Value = -3419.3044

# row 34
$ws.Range("H34").Value = 1815.6666
$ws.Range("I34").Value = 883.12
$ws.Range("J34").Value = 2829.3044
$ws.Range("K34").Value = 883.12
$ws.Range("L34").Value = 2829.3044
$ws.Range("M34").Value = -681.12
$ws.Range("N34").Value = -3233.3044

# row 58
$ws.Range("H58").Value = 2637.5652
$ws.Range("I58").Value = 1044.3529
$ws.Range("J58").Value = 7151.6665
$ws.Range("K58").Value = 1044.3529
$ws.Range("L58").Value = 7151.6665
$ws.Range("M58").Value = -841.3529000000001
$ws.Range("N58").Value = -7557.6665

# row 132
$ws.Range("H132").Value = 1661.8214
$ws.Range("I132").Value = 1383.875
$ws.Range("J132").Value = 2032.4166
$ws.Range("K132").Value = 4151.625
$ws.Range("L132").Value = 6097.2498
$ws.Range("M132").Value = -1621.625
$ws.Range("N132").Value = -11157.2498

# row 134
$ws.Range("H134").Value = 12821462
$ws.Range("I134").Value = 903.9677
$ws.Range("J134").Value = 62501124
$ws.Range("K134").Value = 2711.9031
$ws.Range("L134").Value = 187503372
$ws.Range("M134").Value = -176.9031
$ws.Range("N134").Value = -187508442

# row 136
$ws.Range("H136").Value = 2637.5652
$ws.Range("I136").Value = 1044.3529
$ws.Range("J136").Value = 7151.6665
$ws.Range("K136").Value = 3133.0587
$ws.Range("L136").Value = 21454.9995
$ws.Range("M136").Value = -583.0587000000005
$ws.Range("N136").Value = -26554.9995

$ws = $wb.Worksheets.Item("CUL")
# row 103
$ws.Range("H103").Value = 3179
$ws.Range("I103").Value = 1669.8
$ws.Range("J103").Value = 4122.25
$ws.Range("K103").Value = 5009.4
$ws.Range("L103").Value = 12366.75
$ws.Range("M103").Value = -4130.4
$ws.Range("N103").Value = -14124.75

# row 109
$ws.Range("H109").Value = 1221.7142
$ws.Range("I109").Value = 913.0909
$ws.Range("J109").Value = 2353.3333
$ws.Range("K109").Value = 2739.2727
$ws.Range("L109").Value = 7059.999899999999
$ws.Range("M109").Value = -1699.2727
$ws.Range("N109").Value = -9139.999899999999

# row 112
$ws.Range("H112").Value = 587464
$ws.Range("I112").Value = 980
$ws.Range("J112").Value = 618331.5600000001
$ws.Range("K112").Value = 2940
$ws.Range("L112").Value = 1854994.68
$ws.Range("M112").Value = -1832
$ws.Range("N112").Value = -1857210.68

# row 115
$ws.Range("H115").Value = 1810
$ws.Range("I115").Value = 665
$ws.Range("J115").Value = 4100
$ws.Range("K115").Value = 1995
$ws.Range("L115").Value = 12300
$ws.Range("M115").Value = -820
$ws.Range("N115").Value = -14650

# row 118
$ws.Range("H118").Value = 2967.9092
$ws.Range("J118").Value = 4750
$ws.Range("L118").Value = 14250
$ws.Range("N118").Value = -16736

# row 121
$ws.Range("H121").Value = 639.5
$ws.Range("I121").Value = 353
$ws.Range("J121").Value = 926
$ws.Range("K121").Value = 1059
$ws.Range("L121").Value = 2778
$ws.Range("M121").Value = 251
$ws.Range("N121").Value = -5398

# row 131
$ws.Range("H131").Value = 1006.85187
$ws.Range("J131").Value = 1027.987
$ws.Range("L131").Value = 3083.961
$ws.Range("N131").Value = -13163.961

$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 34793.258
$ws.Range("I132").Value = 2072.5557
$ws.Range("J132").Value = 80098.84
$ws.Range("K132").Value = 6217.6671
$ws.Range("L132").Value = 240296.52
$ws.Range("M132").Value = -3687.6671
$ws.Range("N132").Value = -245356.52

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1043
$ws.Range("I22").Value = 1103.3334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 1103.3334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -808.3334
$ws.Range("N22").Value = -1090

# row 27
$ws.Range("H27").Value = 1043
$ws.Range("I27").Value = 1103.3334
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 1103.3334
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -996.3334
$ws.Range("N27").Value = -714

# row 132
$ws.Range("H132").Value = 548650.7
$ws.Range("I132").Value = 136434.14
$ws.Range("J132").Value = 1431971.9
$ws.Range("K132").Value = 409302.42
$ws.Range("L132").Value = 4295915.699999999
$ws.Range("M132").Value = -406772.42
$ws.Range("N132").Value = -4300975.699999999

$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 4933.375
$ws.Range("I132").Value = 1481.5294
$ws.Range("J132").Value = 8845.467000000001
$ws.Range("K132").Value = 4444.5882
$ws.Range("L132").Value = 26536.401
$ws.Range("M132").Value = -1914.5882
$ws.Range("N132").Value = -31596.401
